# Generate Report for Handoff
#
# A fresh handoff package was generated, so:
#   - the per-locale "Status" column flips from the old handback state to
#     "Ready for handoff"
#   - the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     columns pick up the new generation timestamps
#   - the now-shorter status text means the Status columns can be narrower
#
# "Latest Handback DateTime" columns are untouched -- this run only
# (re)generates the handoff, it doesn't record a new handback.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -----------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-25 06:58:19"

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-25 06:58:13"

# --- de-de sheet ----------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-25 06:58:19"

# --- Column width changes (Status columns narrowed) ------------------
# ColumnWidth snaps to whole-pixel increments of the Normal-style font, so
# this is the closest attainable width to the authored 17.2159881591797.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
